# Weekly update: insert a new "Poroto verde" price record at the top of the
# Vega Modelo de Temuco block (row 125), pushing the existing rows 125:196
# down to 126:197.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("125").Insert()

$ws.Range("A125").Value = 10
$ws.Range("B125").Value = "Vega Modelo de Temuco"
$ws.Range("C125").Value = "La Araucanía"
$ws.Range("D125").Value = 45001
$ws.Range("E125").Value = 9
$ws.Range("F125").Value = 100112031
$ws.Range("G125").Value = "Poroto verde"
$ws.Range("H125").Value = "Brío"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 240
$ws.Range("K125").Value = 1000
$ws.Range("L125").Value = 1200
$ws.Range("M125").Value = 1117
$ws.Range("N125").Value = "$/kilo"
$ws.Range("O125").Value = "Región de La Araucanía"
$ws.Range("P125").Value = 1117
$ws.Range("Q125").Value = 1
$ws.Range("R125").Value = "Hortaliza"
